# Update division problems per commit "Update master to output generated at 1c8df47"
$d = $word.ActiveDocument

# Replacements are ordered so that a newly-written value is never later
# mistaken for an old value still awaiting replacement (e.g. "42÷5=" is
# both a target of one change and the source of another).
$d.Content.Find.Execute("67÷4=", $false, $false, $false, $false, $false, $true, 0, $false, "72÷4=", 1) | Out-Null
$d.Content.Find.Execute("60÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "65÷5=", 1) | Out-Null
$d.Content.Find.Execute("19÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "35÷9=", 1) | Out-Null
$d.Content.Find.Execute("73÷8=", $false, $false, $false, $false, $false, $true, 0, $false, "83÷2=", 1) | Out-Null
$d.Content.Find.Execute("50÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "33÷4=", 1) | Out-Null
$d.Content.Find.Execute("27÷8=", $false, $false, $false, $false, $false, $true, 0, $false, "30÷7=", 1) | Out-Null
$d.Content.Find.Execute("24÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "37÷4=", 1) | Out-Null
$d.Content.Find.Execute("84÷9=", $false, $false, $false, $false, $false, $true, 0, $false, "59÷2=", 1) | Out-Null
$d.Content.Find.Execute("30÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "84÷8=", 1) | Out-Null
$d.Content.Find.Execute("74÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "72÷5=", 1) | Out-Null
$d.Content.Find.Execute("42÷5=", $false, $false, $false, $false, $false, $true, 0, $false, "30÷5=", 1) | Out-Null
$d.Content.Find.Execute("57÷4=", $false, $false, $false, $false, $false, $true, 0, $false, "42÷5=", 1) | Out-Null
$d.Content.Find.Execute("21÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "38÷7=", 1) | Out-Null
$d.Content.Find.Execute("61÷5=", $false, $false, $false, $false, $false, $true, 0, $false, "58÷9=", 1) | Out-Null
$d.Content.Find.Execute("14÷4=", $false, $false, $false, $false, $false, $true, 0, $false, "49÷4=", 1) | Out-Null
$d.Content.Find.Execute("10÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "20÷2=", 1) | Out-Null
$d.Content.Find.Execute("35÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "60÷9=", 1) | Out-Null
$d.Content.Find.Execute("83÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "99÷6=", 1) | Out-Null
$d.Content.Find.Execute("49÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "73÷4=", 1) | Out-Null
$d.Content.Find.Execute("28÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "45÷3=", 1) | Out-Null
$d.Content.Find.Execute("65÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "50÷2=", 1) | Out-Null
$d.Content.Find.Execute("49÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "54÷9=", 1) | Out-Null
$d.Content.Find.Execute("69÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "76÷6=", 1) | Out-Null
$d.Content.Find.Execute("53÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "76÷2=", 1) | Out-Null
$d.Content.Find.Execute("87÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "92÷5=", 1) | Out-Null

Write-Host "All replacements complete"
